# Add EMSL_ID values to the "dilutions_calculations_sheet" sheet, column B,
# for rows 2 through 115 (EUP60985_122 .. EUP60985_235), then switch the
# active sheet / selection to mirror the author's final view state on the
# "print me lab dilution sheet".

$wb = $excel.ActiveWorkbook
$calcSheet = $wb.Worksheets.Item("dilutions_calculations_sheet")
$printSheet = $wb.Worksheets.Item("print me lab dilution sheet")

$startRow = 2
$endRow = 115
$startId = 122

for ($row = $startRow; $row -le $endRow; $row++) {
    $id = $startId + ($row - $startRow)
    $cell = $calcSheet.Cells.Item($row, 2)
    $cell.Value = "EUP60985_" + $id
    $cell.Style = "Normal"
}

# Reflect the scrolled/selected view left in "dilutions_calculations_sheet"
# after entering the data.
$calcSheet.Range("M2").Select()
$excel.ActiveWindow.ScrollColumn = $calcSheet.Range("G1").Column

# Switch to and leave "print me lab dilution sheet" as the active sheet,
# matching the final saved view state.
$printSheet.Activate()
$printSheet.Range("H15").Select()
